$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1930666666666666
$ws.Range("H2").Value = 0.5791999999999999
$ws.Range("I2").Value = 0.01292026122037801
$ws.Range("J2").Value = 0.01292026122037801
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1670496666666667
$ws.Range("N2").Value = 0.501149
$ws.Range("Q2").Value = 0.03225172231111111
$ws.Range("R2").Value = 0.2902655008
$ws.Range("S2").Value = 0.01292026122037801
$ws.Range("T2").Value = 0.01292026122037801

# Row 3
$ws.Range("I3").Value = 0.04457951877603724
$ws.Range("J3").Value = 0.04457951877603725
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1670496666666667
$ws.Range("N3").Value = 0.501149
$ws.Range("Q3").Value = 0.1112799684003333
$ws.Range("R3").Value = 1.001519715603
$ws.Range("S3").Value = 0.04457951877603724
$ws.Range("T3").Value = 0.04457951877603725

# Row 4
$ws.Range("G4").Value = 14.08372266666667
$ws.Range("H4").Value = 42.251168
$ws.Range("I4").Value = 0.9425002200035847
$ws.Range("J4").Value = 0.9425002200035848
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1670496666666667
$ws.Range("N4").Value = 0.501149
$ws.Range("Q4").Value = 2.352681176892444
$ws.Range("R4").Value = 21.174130592032
$ws.Range("S4").Value = 0.9425002200035847
$ws.Range("T4").Value = 0.9425002200035848
